$d = $word.ActiveDocument

# wdReplaceAll = 2 (last arg), MatchCase=$true, MatchWholeWord=$false, MatchWildcards=$false,
# MatchSoundsLike=$false, MatchAllWordForms=$false, Forward=$true, Wrap=1 (wdFindContinue),
# Format=$false
$wdFindContinue = 1
$wdReplaceAll = 2

# --- Title: "Relátorio de progresso 22 – 28 Maio" -> "Relátorio de progresso 29 – 5 Junho" ---
$d.Content.Find.Execute(
    "Relátorio de progresso 22 – 28 Maio", $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "Relátorio de progresso 29 – 5 Junho", $wdReplaceAll)

# Move the "_GoBack" bookmark from the end of bullet 1 to the end of the (new) title text.
# Adding a bookmark with the existing name "_GoBack" relocates it (removing the old one),
# mirroring Word's own "last edit location" bookmark behaviour.
$p1 = $d.Paragraphs(1)
$titleEnd = $p1.Range.End - 1
$marker = $d.Range($titleEnd, $titleEnd)
$marker.InsertAfter("~")
$markerRng = $d.Range($titleEnd, $titleEnd + 1)
$d.Bookmarks.Add("_GoBack", $markerRng)
$d.Range($titleEnd, $titleEnd + 1).Delete()

# --- Bullet 1: Setup... -> Leitura sobre sincronização de Broadcast receivers em Android ---
$d.Content.Find.Execute(
    "Setup de todas as ferramentas usadas no desenvolvimento.", $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "Leitura sobre sincronização de Broadcast receivers em Android", $wdReplaceAll)

# --- Bullet 2: Implementação dos pedidos do servidor... -> Leitura sobre Multipart ---
$d.Content.Find.Execute(
    "Implementação dos pedidos do servidor aos servidores da google, para a realização de push de informação aos dispositivos móveis.",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "Leitura sobre Multipart", $wdReplaceAll)

# --- Bullet 3: Implementação da lógica da aplicação Android... -> Implementação de um endpoint... ---
$d.Content.Find.Execute(
    "Implementação da lógica da aplicação Android quando recebe um pedido dos servidores da Google(pedido firebase).",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "Implementação de um endpoint que permite enviar ficheiros para o servidor", $wdReplaceAll)

# --- Bullet 4: Implementação do push de informação... -> Implementação da deteção de imagens copiadas em Windows (trailing space) ---
$d.Content.Find.Execute(
    "Implementação do push de informação para o servidor na aplicação Windows.",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "Implementação da deteção de imagens copiadas em Windows ", $wdReplaceAll)

# --- Bullet 5: Inicio da implementação de thread pool... -> Implementação de metodo que realiza fetch... ---
$d.Content.Find.Execute(
    "Inicio da implementação de thread pool da informação contida no servidor na aplicação Windows.",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "Implementação de metodo que realiza fetch da informação na aplicação Windows", $wdReplaceAll)
